$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the text (string) storage type for Price/Volume columns —
# the source workbook stores these as plain text (e.g. "277.23", "0.9680"),
# not numbers, so force text format before writing to avoid Excel
# auto-converting numeric-looking strings into numeric cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "20.586.35"
$ws.Range("E2").Value = "  +2.12%  "

$ws.Range("D3").Value = "1.479.93"
$ws.Range("E3").Value = "  +3.62%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "0.9680"
$ws.Range("E5").Value = "  -2.85%  "

$ws.Range("D6").Value = "277.23"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "0.3659"
$ws.Range("E7").Value = "  -1.14%  "

$ws.Range("D8").Value = "0.3066"
$ws.Range("E8").Value = "  -2.75%  "

$ws.Range("D9").Value = "40.65"
$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("D10").Value = "1.063"
$ws.Range("E10").Value = "  +0.52%  "

$ws.Range("D11").Value = "0.06650"
$ws.Range("E11").Value = "  +0.99%  "

$ws.Range("D12").Value = "0.9970"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").Value = "5.497"
$ws.Range("E13").Value = "  -1.00%  "

$ws.Range("D14").Value = "18.21"
$ws.Range("E14").Value = "  -0.03%  "

$ws.Range("D15").Value = "6.195"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").Value = "0.00001035"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").Value = "1.479.46"
$ws.Range("E17").Value = "  +3.30%  "

$ws.Range("D18").Value = "0.9693"
$ws.Range("E18").Value = "  -2.69%  "

$ws.Range("D19").Value = "0.05923"
$ws.Range("E19").Value = "  +3.06%  "

$ws.Range("D20").Value = "69.61"
$ws.Range("E20").Value = "  -2.90%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "14.62"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.460"
$ws.Range("E22").Value = "  -2.96%  "

$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("D24").Value = "2.258"
$ws.Range("E24").Value = "  +0.95%  "

$ws.Range("D25").Value = "20.614.26"
$ws.Range("E25").Value = "  +1.96%  "

$ws.Range("D26").Value = "141.00"
$ws.Range("E26").Value = "  +4.08%  "

$ws.Range("D27").Value = "2.147"
$ws.Range("E27").Value = "  -7.20%  "

$ws.Range("D28").Value = "17.33"
$ws.Range("E28").Value = "  -0.78%  "

$ws.Range("D29").Value = "1.635.44"
$ws.Range("E29").Value = "  +2.63%  "

$ws.Range("D30").Value = "114.09"
$ws.Range("E30").Value = "  +2.02%  "

$ws.Range("D31").Value = "3.938"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").Value = "0.8175"
$ws.Range("E32").Value = "  -2.76%  "

$ws.Range("D33").Value = "4.983"
$ws.Range("E33").Value = "  -6.26%  "

$ws.Range("D34").Value = "0.07926"
$ws.Range("E34").Value = "  +1.54%  "

$ws.Range("D35").Value = "1.553"
$ws.Range("E35").Value = "  +3.74%  "

$ws.Range("D36").Value = "1.215"
$ws.Range("E36").Value = "  +9.33%  "

$ws.Range("D37").Value = "0.05846"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("D38").Value = "4.740"
$ws.Range("E38").Value = "  -3.83%  "

$ws.Range("D39").Value = "0.9689"
$ws.Range("E39").Value = "  -2.78%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02044"
$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "10.49"
$ws.Range("E41").Value = "  -2.63%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "7.660"
$ws.Range("E42").Value = "  -2.19%  "

$ws.Range("D43").Value = "0.1884"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").Value = "0.5310"
$ws.Range("E44").Value = "  -1.18%  "

$ws.Range("D45").Value = "3.510"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("D46").Value = "12.20"
$ws.Range("E46").Value = "  -2.15%  "

$ws.Range("D47").Value = "118.29"
$ws.Range("E47").Value = "  -1.58%  "

$ws.Range("D48").Value = "0.5212"
$ws.Range("E48").Value = "  -1.28%  "

$ws.Range("D49").Value = "1.804"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").Value = "0.06469"
$ws.Range("E50").Value = "  +2.89%  "

$ws.Range("D51").Value = "0.9931"
$ws.Range("E51").Value = "  -0.65%  "
